# Re-upload of the two SharePoint/"document management" custom XML parts that
# ship with this AICTE PPT Template.
#
# The package carries a pair of linked customXml parts:
#   - one holding the "FormTemplates" (mso-contentType) SharePoint forms
#     declaration, tied to schema http://schemas.microsoft.com/sharepoint/v3/contenttype/forms
#     (itemID {927BD4C1-B6B1-4715-ABF9-E660A51A4EA0})
#   - one holding the "p:properties / documentManagement / _activity" part
#     (itemID {8D289AE2-D2AE-49D1-AFAC-3A79F6794255})
#
# The upload re-ordered which numbered part (item1.xml/item2.xml) holds which
# payload: the FormTemplates part now comes first and the properties part
# second. Reproduce that by removing both parts and re-adding them, through
# the CustomXMLParts collection, in the desired order so PowerPoint's save
# renumbers customXml/item1.xml .. item2.xml (and their itemProps*.xml
# companions) to match.

$p = $ppt.ActivePresentation
$parts = $p.CustomXMLParts

$formsId = "{927BD4C1-B6B1-4715-ABF9-E660A51A4EA0}"
$propsId = "{8D289AE2-D2AE-49D1-AFAC-3A79F6794255}"

$formsXml = "<?mso-contentType ?>`r`n<FormTemplates xmlns=`"http://schemas.microsoft.com/sharepoint/v3/contenttype/forms`">`r`n  <Display>DocumentLibraryForm</Display>`r`n  <Edit>DocumentLibraryForm</Edit>`r`n  <New>DocumentLibraryForm</New>`r`n</FormTemplates>"

$propsXml = "<p:properties xmlns:p=`"http://schemas.microsoft.com/office/2006/metadata/properties`" xmlns:xsi=`"http://www.w3.org/2001/XMLSchema-instance`" xmlns:pc=`"http://schemas.microsoft.com/office/infopath/2007/PartnerControls`">`r`n  <documentManagement>`r`n    <_activity xmlns=`"b30265f8-c5e2-4918-b4a1-b977299ca3e2`" xsi:nil=`"true`"/>`r`n  </documentManagement>`r`n</p:properties>"

function Remove-PartById($collection, $id) {
    try {
        $existing = $collection.SelectByID($id)
        if ($existing -ne $null) {
            $existing.Delete()
        }
    } catch {
        # Part not present / not addressable by id - nothing to remove.
    }
}

# Drop the previous copies of both linked parts ...
Remove-PartById $parts $formsId
Remove-PartById $parts $propsId

# ... and re-add them with the FormTemplates part first, so it lands in
# customXml/item1.xml (companion customXml/itemProps1.xml), and the
# properties part second, landing in customXml/item2.xml (companion
# customXml/itemProps2.xml) - matching the re-uploaded file layout.
$parts.Add($formsXml) | Out-Null
$parts.Add($propsXml) | Out-Null
